# Auto-generated: apply 2026-01-14 violent crime data updates
# Each entry: worksheet name -> list of (cell reference, new numeric value)
$wb = $excel.ActiveWorkbook

$updates = @{}

$updates['Citywide Totals'] = @(
    @{ Cell = 'M2'; Value = 205 }
    @{ Cell = 'L3'; Value = 7112 }
    @{ Cell = 'M3'; Value = 229 }
    @{ Cell = 'L4'; Value = 1787 }
    @{ Cell = 'M4'; Value = 51 }
    @{ Cell = 'M6'; Value = 175 }
    @{ Cell = 'M7'; Value = 671 }
)

$updates['Norwood Park'] = @(
    @{ Cell = 'M3'; Value = 2 }
    @{ Cell = 'M7'; Value = 3 }
)

$updates['Logan Square'] = @(
    @{ Cell = 'M3'; Value = 3 }
    @{ Cell = 'M7'; Value = 6 }
)

$updates['Austin'] = @(
    @{ Cell = 'L4'; Value = 103 }
    @{ Cell = 'M6'; Value = 13 }
    @{ Cell = 'L7'; Value = 1434 }
    @{ Cell = 'M7'; Value = 39 }
)

$updates['West Pullman'] = @(
    @{ Cell = 'M6'; Value = 2 }
    @{ Cell = 'M7'; Value = 7 }
)

$updates['Grand Crossing'] = @(
    @{ Cell = 'M3'; Value = 11 }
    @{ Cell = 'M4'; Value = 5 }
    @{ Cell = 'M7'; Value = 33 }
)

$updates['Woodlawn'] = @(
    @{ Cell = 'M2'; Value = 7 }
    @{ Cell = 'M7'; Value = 22 }
)

$updates['By Neighborhood'] = @(
    @{ Cell = 'M6'; Value = 6 }
    @{ Cell = 'M7'; Value = 21 }
    @{ Cell = 'L8'; Value = 1434 }
    @{ Cell = 'M8'; Value = 39 }
    @{ Cell = 'M9'; Value = 8 }
    @{ Cell = 'M10'; Value = 3 }
    @{ Cell = 'M11'; Value = 9 }
    @{ Cell = 'L16'; Value = 48 }
    @{ Cell = 'M18'; Value = 5 }
    @{ Cell = 'M19'; Value = 23 }
    @{ Cell = 'M20'; Value = 28 }
    @{ Cell = 'L27'; Value = 189 }
    @{ Cell = 'L29'; Value = 1213 }
    @{ Cell = 'M29'; Value = 33 }
    @{ Cell = 'M31'; Value = 11 }
    @{ Cell = 'M36'; Value = 7 }
    @{ Cell = 'M37'; Value = 33 }
    @{ Cell = 'M41'; Value = 5 }
    @{ Cell = 'M42'; Value = 18 }
    @{ Cell = 'M44'; Value = 4 }
    @{ Cell = 'M47'; Value = 6 }
    @{ Cell = 'M51'; Value = 9 }
    @{ Cell = 'L52'; Value = 457 }
    @{ Cell = 'M53'; Value = 6 }
    @{ Cell = 'M60'; Value = 7 }
    @{ Cell = 'L63'; Value = 65 }
    @{ Cell = 'M64'; Value = 11 }
    @{ Cell = 'M69'; Value = 3 }
    @{ Cell = 'M79'; Value = 16 }
    @{ Cell = 'M85'; Value = 37 }
    @{ Cell = 'L89'; Value = 291 }
    @{ Cell = 'M89'; Value = 6 }
    @{ Cell = 'M91'; Value = 6 }
    @{ Cell = 'M94'; Value = 8 }
    @{ Cell = 'M95'; Value = 7 }
    @{ Cell = 'M96'; Value = 8 }
    @{ Cell = 'M97'; Value = 8 }
    @{ Cell = 'M99'; Value = 22 }
    @{ Cell = 'M101'; Value = 671 }
)

$updates['Gage Park'] = @(
    @{ Cell = 'M2'; Value = 5 }
    @{ Cell = 'M7'; Value = 11 }
)

$updates['Englewood'] = @(
    @{ Cell = 'M2'; Value = 13 }
    @{ Cell = 'L4'; Value = 66 }
    @{ Cell = 'L7'; Value = 1213 }
    @{ Cell = 'M7'; Value = 33 }
)

$updates['Chatham'] = @(
    @{ Cell = 'M2'; Value = 5 }
    @{ Cell = 'M6'; Value = 7 }
    @{ Cell = 'M7'; Value = 23 }
)

$updates['Irving Park'] = @(
    @{ Cell = 'M4'; Value = 1 }
    @{ Cell = 'M7'; Value = 4 }
)

$updates['Ashburn'] = @(
    @{ Cell = 'M6'; Value = 3 }
    @{ Cell = 'M7'; Value = 6 }
)

$updates['Hermosa'] = @(
    @{ Cell = 'M2'; Value = 2 }
    @{ Cell = 'M7'; Value = 5 }
)

$updates['Humboldt Park'] = @(
    @{ Cell = 'M3'; Value = 9 }
    @{ Cell = 'M7'; Value = 18 }
)

$updates['Avondale'] = @(
    @{ Cell = 'M2'; Value = 1 }
    @{ Cell = 'M7'; Value = 3 }
)

$updates['West Ridge'] = @(
    @{ Cell = 'M2'; Value = 2 }
    @{ Cell = 'M7'; Value = 8 }
)

$updates['Washington Park'] = @(
    @{ Cell = 'M3'; Value = 4 }
    @{ Cell = 'M7'; Value = 6 }
)

$updates['Roseland'] = @(
    @{ Cell = 'M3'; Value = 4 }
    @{ Cell = 'M7'; Value = 16 }
)

$updates['Near South Side'] = @(
    @{ Cell = 'M2'; Value = 2 }
    @{ Cell = 'M7'; Value = 11 }
)

$updates['Chicago Lawn'] = @(
    @{ Cell = 'M2'; Value = 10 }
    @{ Cell = 'M4'; Value = 3 }
    @{ Cell = 'M7'; Value = 28 }
)

$updates['Calumet Heights'] = @(
    @{ Cell = 'M3'; Value = 3 }
    @{ Cell = 'M7'; Value = 5 }
)

$updates['Grand Boulevard'] = @(
    @{ Cell = 'M2'; Value = 4 }
    @{ Cell = 'M3'; Value = 3 }
    @{ Cell = 'M7'; Value = 7 }
)

$updates['Auburn Gresham'] = @(
    @{ Cell = 'M3'; Value = 9 }
    @{ Cell = 'M7'; Value = 21 }
)

$updates['West Loop'] = @(
    @{ Cell = 'M6'; Value = 5 }
    @{ Cell = 'M7'; Value = 8 }
)

$updates['Kenwood'] = @(
    @{ Cell = 'M3'; Value = 4 }
    @{ Cell = 'M7'; Value = 6 }
)

$updates['Belmont Cragin'] = @(
    @{ Cell = 'M6'; Value = 3 }
    @{ Cell = 'M7'; Value = 9 }
)

$updates['Avalon Park'] = @(
    @{ Cell = 'M2'; Value = 2 }
    @{ Cell = 'M7'; Value = 8 }
)

$updates['West Town'] = @(
    @{ Cell = 'M4'; Value = 2 }
    @{ Cell = 'M7'; Value = 8 }
)

$updates['Uptown'] = @(
    @{ Cell = 'M2'; Value = 2 }
    @{ Cell = 'M3'; Value = 2 }
    @{ Cell = 'L6'; Value = 82 }
    @{ Cell = 'L7'; Value = 291 }
    @{ Cell = 'M7'; Value = 6 }
)

$updates['Edgewater'] = @(
    @{ Cell = 'L4'; Value = 26 }
    @{ Cell = 'L7'; Value = 189 }
)

$updates['Little Italy, UIC'] = @(
    @{ Cell = 'M6'; Value = 2 }
    @{ Cell = 'M7'; Value = 9 }
)

$updates['Morgan Park'] = @(
    @{ Cell = 'M2'; Value = 3 }
    @{ Cell = 'M3'; Value = 2 }
    @{ Cell = 'M6'; Value = 1 }
    @{ Cell = 'M7'; Value = 7 }
)

$updates['South Shore'] = @(
    @{ Cell = 'M2'; Value = 8 }
    @{ Cell = 'M3'; Value = 18 }
    @{ Cell = 'M7'; Value = 37 }
)

$updates['Little Village'] = @(
    @{ Cell = 'L4'; Value = 29 }
    @{ Cell = 'L6'; Value = 128 }
    @{ Cell = 'L7'; Value = 457 }
)

$updates['Bucktown'] = @(
    @{ Cell = 'L6'; Value = 29 }
    @{ Cell = 'L7'; Value = 48 }
)

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}
